$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared strings used in the data table.
# Column K (rows 2-29) held "sig_eta_Z" -> rename to "sig_Z_eta"
# Column L header (L1) held "ixsec" -> rename to "normalization"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 11).Value = "sig_Z_eta"
}
$ws.Cells.Item(1, 12).Value = "normalization"

# Update the L column (rows 2-29) data values from 1 to 256.60000000000002
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 12).Value = 256.60000000000002
}

# Widen column L (12) to fit the new header text (resulting stored width = 14)
$ws.Columns.Item(12).ColumnWidth = 13.166666666666666

# Update the active selection to L31
$ws.Range("L31").Select()
